# Apply the daily cryptos data refresh: updates Price (D) and Volume(1h) (E)
# values for rows 2-51, plus the three-row insertion of BabyDogeCoin ahead of
# Stellar/FirstDigitalUSD/ApeXProtocol (rows 45-48 shift down by one coin).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to stay plain text so values like "1.00" / "431.10" /
# "68.375.46" are not silently reinterpreted as numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.375.46'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.904.63'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '485.34'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("E6").Value = '  +1.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.44%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.742'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.68%  '
$ws.Range("E10").Value = '  +8.08%  '
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.03'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("E13").Value = '  -0.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.519.42'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.937.70'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.23'
$ws.Range("D16").ClearFormats()
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.98'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.449.54'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '431.10'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("E22").Value = '  +6.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.77'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.38'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +21.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '89.28'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.27%  '
$ws.Range("E26").Value = '  +4.08%  '
$ws.Range("E27").Value = '  -5.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.35'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.70'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '716.78'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("E32").Value = '  +0.46%  '
$ws.Range("E33").Value = '  +2.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0895'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '61.65'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +5.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.08'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +8.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.80'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.405'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +19.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.148'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.79%  '
$ws.Range("E41").Value = '  +4.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.99'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +9.67%  '
$ws.Range("E43").Value = '  +3.76%  '
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₆0369'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +26.62%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.142'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.32'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +5.34%  '
$ws.Range("E49").Value = '  -1.10%  '
$ws.Range("E50").Value = '  -2.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.33'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.76%  '
